$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BOLIVIA - DIVISION PROFESIONAL row (row 3); rows below shift up.
$ws.Rows(3).Delete()

# --- Row 2 (Argentina - Ind. Rivadavia vs River Plate) odds tweaks ---
$ws.Range("I2").Value = 1.44
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73
$ws.Range("W2").Value = 15
$ws.Range("Y2").Value = 23
$ws.Range("AP2").Value = 41

# --- Row 3 (was row 4: Colombia - Junior vs America De Cali) odds tweaks ---
$ws.Range("G3").Value = 2.25
$ws.Range("I3").Value = 3.5
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("AF3").Value = 51
$ws.Range("AV3").Value = 5

# --- Row 4 (was row 5: Mexico - Guadalajara Chivas vs Atlas) odds tweaks ---
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 2.05
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 7.5
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6
$ws.Range("W4").Value = 5.5
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 7.5
$ws.Range("AH4").Value = 10
$ws.Range("AO4").Value = 9.5
$ws.Range("AR4").Value = 67
$ws.Range("AT4").Value = 9.5
$ws.Range("AY4").Value = 126
$ws.Range("BA4").Value = 401
